# Updated under ice read length
# Forward read length (BY) and Reverse read length (BZ) for rows 12-32 and
# 34-35 change from 125 to 124 (row 33 is a blank spacer row and is skipped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(12..32) + @(34,35)

foreach ($r in $rows) {
    $ws.Range("BY$r").Value = 124
    $ws.Range("BZ$r").Value = 124
}

# Leave the selection where the editing left off, matching the saved file.
$ws.Range("CA33").Select()
